$wb = $excel.ActiveWorkbook

# --- Create the new "2022-Q4" worksheet with fund holdings data ---
$q4 = $wb.Worksheets.Add()

# Header row
$q4.Cells.Item(1,2).Value = "基金代码"
$q4.Cells.Item(1,3).Value = "基金名称"
$q4.Cells.Item(1,4).Value = "基金规模"
$q4.Cells.Item(1,5).Value = "股票总仓位"
$q4.Cells.Item(1,6).Value = "仓位占比"
$q4.Cells.Item(1,7).Value = "持有市值(亿元)"
$q4.Cells.Item(1,8).Value = "仓位排名"

# Data rows
$q4.Cells.Item(2,1).Value = 0
$q4.Cells.Item(2,2).NumberFormat = "@"
$q4.Cells.Item(2,2).Value = "002770"
$q4.Cells.Item(2,3).Value = "安信新回报灵活配置混合A"
$q4.Cells.Item(2,4).NumberFormat = "@"
$q4.Cells.Item(2,4).Value = "2.74"
$q4.Cells.Item(2,5).NumberFormat = "@"
$q4.Cells.Item(2,5).Value = "81.00"
$q4.Cells.Item(2,6).NumberFormat = "@"
$q4.Cells.Item(2,6).Value = "4.29"
$q4.Cells.Item(2,7).NumberFormat = "@"
$q4.Cells.Item(2,7).Value = "0.1175"
$q4.Cells.Item(2,8).Value = 7

$q4.Cells.Item(3,1).Value = 1
$q4.Cells.Item(3,2).NumberFormat = "@"
$q4.Cells.Item(3,2).Value = "005014"
$q4.Cells.Item(3,3).Value = "泰康景泰回报混合A"
$q4.Cells.Item(3,4).NumberFormat = "@"
$q4.Cells.Item(3,4).Value = "8.86"
$q4.Cells.Item(3,5).NumberFormat = "@"
$q4.Cells.Item(3,5).Value = "32.91"
$q4.Cells.Item(3,6).NumberFormat = "@"
$q4.Cells.Item(3,6).Value = "1.20"
$q4.Cells.Item(3,7).NumberFormat = "@"
$q4.Cells.Item(3,7).Value = "0.1063"
$q4.Cells.Item(3,8).Value = 9

$q4.Cells.Item(4,1).Value = 2
$q4.Cells.Item(4,2).NumberFormat = "@"
$q4.Cells.Item(4,2).Value = "001449"
$q4.Cells.Item(4,3).Value = "华商双驱优选灵活配置混合"
$q4.Cells.Item(4,4).NumberFormat = "@"
$q4.Cells.Item(4,4).Value = "2.26"
$q4.Cells.Item(4,5).NumberFormat = "@"
$q4.Cells.Item(4,5).Value = "77.71"
$q4.Cells.Item(4,6).NumberFormat = "@"
$q4.Cells.Item(4,6).Value = "4.07"
$q4.Cells.Item(4,7).NumberFormat = "@"
$q4.Cells.Item(4,7).Value = "0.0920"
$q4.Cells.Item(4,8).Value = 5

$q4.Cells.Item(5,1).Value = 3
$q4.Cells.Item(5,2).NumberFormat = "@"
$q4.Cells.Item(5,2).Value = "002771"
$q4.Cells.Item(5,3).Value = "安信新回报灵活配置混合C"
$q4.Cells.Item(5,4).NumberFormat = "@"
$q4.Cells.Item(5,4).Value = "2.10"
$q4.Cells.Item(5,5).NumberFormat = "@"
$q4.Cells.Item(5,5).Value = "81.00"
$q4.Cells.Item(5,6).NumberFormat = "@"
$q4.Cells.Item(5,6).Value = "4.29"
$q4.Cells.Item(5,7).NumberFormat = "@"
$q4.Cells.Item(5,7).Value = "0.0901"
$q4.Cells.Item(5,8).Value = 7

$q4.Cells.Item(6,1).Value = 4
$q4.Cells.Item(6,2).NumberFormat = "@"
$q4.Cells.Item(6,2).Value = "010403"
$q4.Cells.Item(6,3).Value = "华商景气优选混合"
$q4.Cells.Item(6,4).NumberFormat = "@"
$q4.Cells.Item(6,4).Value = "0.58"
$q4.Cells.Item(6,5).NumberFormat = "@"
$q4.Cells.Item(6,5).Value = "77.20"
$q4.Cells.Item(6,6).NumberFormat = "@"
$q4.Cells.Item(6,6).Value = "3.88"
$q4.Cells.Item(6,7).NumberFormat = "@"
$q4.Cells.Item(6,7).Value = "0.0225"
$q4.Cells.Item(6,8).Value = 6

$q4.Cells.Item(7,1).Value = 5
$q4.Cells.Item(7,2).NumberFormat = "@"
$q4.Cells.Item(7,2).Value = "014627"
$q4.Cells.Item(7,3).Value = "财通多策略福瑞混合（LOF）C"
$q4.Cells.Item(7,4).NumberFormat = "@"
$q4.Cells.Item(7,4).Value = "1.09"
$q4.Cells.Item(7,5).NumberFormat = "@"
$q4.Cells.Item(7,5).Value = "62.26"
$q4.Cells.Item(7,6).NumberFormat = "@"
$q4.Cells.Item(7,6).Value = "1.15"
$q4.Cells.Item(7,7).NumberFormat = "@"
$q4.Cells.Item(7,7).Value = "0.0125"
$q4.Cells.Item(7,8).Value = 7

$q4.Cells.Item(8,1).Value = 6
$q4.Cells.Item(8,2).NumberFormat = "@"
$q4.Cells.Item(8,2).Value = "501028"
$q4.Cells.Item(8,3).Value = "财通多策略福瑞混合（LOF）A"
$q4.Cells.Item(8,4).NumberFormat = "@"
$q4.Cells.Item(8,4).Value = "0.81"
$q4.Cells.Item(8,5).NumberFormat = "@"
$q4.Cells.Item(8,5).Value = "62.26"
$q4.Cells.Item(8,6).NumberFormat = "@"
$q4.Cells.Item(8,6).Value = "1.15"
$q4.Cells.Item(8,7).NumberFormat = "@"
$q4.Cells.Item(8,7).Value = "0.0093"
$q4.Cells.Item(8,8).Value = 7

$q4.Cells.Item(9,1).Value = 7
$q4.Cells.Item(9,2).NumberFormat = "@"
$q4.Cells.Item(9,2).Value = "005015"
$q4.Cells.Item(9,3).Value = "泰康景泰回报混合C"
$q4.Cells.Item(9,4).NumberFormat = "@"
$q4.Cells.Item(9,4).Value = "0.37"
$q4.Cells.Item(9,5).NumberFormat = "@"
$q4.Cells.Item(9,5).Value = "32.91"
$q4.Cells.Item(9,6).NumberFormat = "@"
$q4.Cells.Item(9,6).Value = "1.20"
$q4.Cells.Item(9,7).NumberFormat = "@"
$q4.Cells.Item(9,7).Value = "0.0044"
$q4.Cells.Item(9,8).Value = 9

$q4.Cells.Item(10,1).Value = 8
$q4.Cells.Item(10,2).NumberFormat = "@"
$q4.Cells.Item(10,2).Value = "008300"
$q4.Cells.Item(10,3).Value = "人保量化锐进混合A"
$q4.Cells.Item(10,4).NumberFormat = "@"
$q4.Cells.Item(10,4).Value = "0.08"
$q4.Cells.Item(10,5).NumberFormat = "@"
$q4.Cells.Item(10,5).Value = "90.63"
$q4.Cells.Item(10,6).NumberFormat = "@"
$q4.Cells.Item(10,6).Value = "4.92"
$q4.Cells.Item(10,7).NumberFormat = "@"
$q4.Cells.Item(10,7).Value = "0.0039"
$q4.Cells.Item(10,8).Value = 2

$q4.Cells.Item(11,1).Value = 9
$q4.Cells.Item(11,2).NumberFormat = "@"
$q4.Cells.Item(11,2).Value = "006226"
$q4.Cells.Item(11,3).Value = "人保量化基本面混合C"
$q4.Cells.Item(11,4).NumberFormat = "@"
$q4.Cells.Item(11,4).Value = "0.06"
$q4.Cells.Item(11,5).NumberFormat = "@"
$q4.Cells.Item(11,5).Value = "79.86"
$q4.Cells.Item(11,6).NumberFormat = "@"
$q4.Cells.Item(11,6).Value = "3.81"
$q4.Cells.Item(11,7).NumberFormat = "@"
$q4.Cells.Item(11,7).Value = "0.0023"
$q4.Cells.Item(11,8).Value = 3

$q4.Cells.Item(12,1).Value = 10
$q4.Cells.Item(12,2).NumberFormat = "@"
$q4.Cells.Item(12,2).Value = "008301"
$q4.Cells.Item(12,3).Value = "人保量化锐进混合C"
$q4.Cells.Item(12,4).NumberFormat = "@"
$q4.Cells.Item(12,4).Value = "0.04"
$q4.Cells.Item(12,5).NumberFormat = "@"
$q4.Cells.Item(12,5).Value = "90.63"
$q4.Cells.Item(12,6).NumberFormat = "@"
$q4.Cells.Item(12,6).Value = "4.92"
$q4.Cells.Item(12,7).NumberFormat = "@"
$q4.Cells.Item(12,7).Value = "0.0020"
$q4.Cells.Item(12,8).Value = 2

$q4.Cells.Item(13,1).Value = 11
$q4.Cells.Item(13,2).NumberFormat = "@"
$q4.Cells.Item(13,2).Value = "006225"
$q4.Cells.Item(13,3).Value = "人保量化基本面混合A"
$q4.Cells.Item(13,4).NumberFormat = "@"
$q4.Cells.Item(13,4).Value = "0.01"
$q4.Cells.Item(13,5).NumberFormat = "@"
$q4.Cells.Item(13,5).Value = "79.86"
$q4.Cells.Item(13,6).NumberFormat = "@"
$q4.Cells.Item(13,6).Value = "3.81"
$q4.Cells.Item(13,7).NumberFormat = "@"
$q4.Cells.Item(13,7).Value = "0.0004"
$q4.Cells.Item(13,8).Value = 3

# Name + position it right after "总计" (as sheet 2, before "2022-Q3")
$q4.Name = "2022-Q4"
$q4.Move($wb.Worksheets.Item(2))

# --- Update the "总计" (summary) sheet: insert a row for 2022-Q4 ---
$summary = $wb.Worksheets.Item(1)
$summary.Rows.Item(2).Insert()
$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q4"
$summary.Cells.Item(2,3).Value = 12
$summary.Cells.Item(2,4).Value = 0.46

# Renumber the index column for the rows that shifted down
$summary.Cells.Item(3,1).Value = 1
$summary.Cells.Item(4,1).Value = 2
$summary.Cells.Item(5,1).Value = 3
$summary.Cells.Item(6,1).Value = 4

Write-Output "2022-Q4 sheet added and 总计 updated"
